$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 46, shifting existing rows 46:61 down to 47:62.
$ws.Rows("46:46").Insert()

# Populate the newly inserted row 46 with the new record.
$ws.Range("A46").Value = 1
$ws.Range("B46").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C46").Value = "Arica y Parinacota"
$ws.Range("D46").Value = 44924
$ws.Range("E46").Value = 15
$ws.Range("F46").Value = "Fruta"
$ws.Range("G46").Value = 100103
$ws.Range("H46").Value = "Frutos de hueso (carozo)"
$ws.Range("I46").Value = 100103004
$ws.Range("J46").Value = "Durazno"
$ws.Range("K46").Value = "Springtime"
$ws.Range("L46").Value = "Primera"
$ws.Range("M46").Value = 350
$ws.Range("N46").Value = 23000
$ws.Range("O46").Value = 25000
$ws.Range("P46").Value = 23857
$ws.Range("Q46").Value = "$/bandeja 18 kilos granel"
$ws.Range("R46").Value = "Región de O'Higgins"
$ws.Range("S46").Value = 1325
$ws.Range("T46").Value = 18
